# Adds 12 new Friday, Jan 13 departure rows (rows 130-141) to the
# "Main Data" sheet of the POZ_Departures workbook, mirroring the
# existing table layout: NUMBER, DATE, TIME, FLIGHT, TO, SHORT,
# AIRLINE, MODEL, AIRCFAT ID, STATUS, (blank), DIFFERENCE, (blank).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{Row=130; A=129; B="Friday, Jan 13"; C="2:25 PM"; D="SK1756"; E="Copenhagen"; F="(CPH)"; G="SAS "; H="CRJ9"; I="(EI-FPU)"; J="2:34 PM"; L="0 hours, 9 minutes"}
    @{Row=131; A=130; B="Friday, Jan 13"; C="2:45 PM"; D="LO3944"; E="Warsaw"; F="(WAW)"; G="LOT "; H="E170"; I="(SP-LDF)"; J="2:57 PM"; L="0 hours, 12 minutes"}
    @{Row=132; A=131; B="Friday, Jan 13"; C="3:05 PM"; D="KL1274"; E="Amsterdam"; F="(AMS)"; G="KLM "; H="E75L"; I="(PH-EXS)"; J="3:08 PM"; L="0 hours, 3 minutes"}
    @{Row=133; A=132; B="Friday, Jan 13"; C="4:05 PM"; D="FR7948"; E="Bristol"; F="(BRS)"; G="Ryanair "; H="B738"; I="(SP-RSM)"; J="4:08 PM"; L="0 hours, 3 minutes"}
    @{Row=134; A=133; B="Friday, Jan 13"; C="4:30 PM"; D="W91901"; E="London"; F="(LTN)"; G="Wizz Air "; H="A320"; I="(G-WUKF)"; J="4:28 PM"; L="0 hours, -2 minutes"}
    @{Row=135; A=134; B="Friday, Jan 13"; C="5:30 PM"; D="FR1975"; E="Dublin"; F="(DUB)"; G="Ryanair "; H="B738"; I="(SP-RKR)"; J="5:37 PM"; L="0 hours, 7 minutes"}
    @{Row=136; A=135; B="Friday, Jan 13"; C="5:55 PM"; D="FR7889"; E="Malta"; F="(MLA)"; G="Ryanair "; H="B738"; I="(SP-RSX)"; J="6:00 PM"; L="0 hours, 5 minutes"}
    @{Row=137; A=136; B="Friday, Jan 13"; C="6:20 PM"; D="LO3948"; E="Warsaw"; F="(WAW)"; G="LOT "; H="E170"; I="(SP-LDI)"; J="6:17 PM"; L="0 hours, -3 minutes"}
    @{Row=138; A=137; B="Friday, Jan 13"; C="9:40 PM"; D="P81987"; E="Cologne"; F="(CGN)"; G="SprintAir "; H="AT73"; I="(SP-SPD)"; J="9:43 PM"; L="0 hours, 3 minutes"}
    @{Row=139; A=138; B="Friday, Jan 13"; C="9:45 PM"; D="FR7679"; E="Stockholm"; F="(ARN)"; G="Ryanair "; H="B38M"; I="(9H-VUJ)"; J="10:27 PM"; L="0 hours, 42 minutes"}
    @{Row=140; A=139; B="Friday, Jan 13"; C="10:05 PM"; D="FR8325"; E="London"; F="(STN)"; G="Ryanair "; H="B38M"; I="(EI-HEY)"; J="10:33 PM"; L="0 hours, 28 minutes"}
    @{Row=141; A=140; B="Friday, Jan 13"; C="10:10 PM"; D="FR3594"; E="Milan"; F="(BGY)"; G="Malta Air "; H="B38M"; I="(9H-VUB)"; J="10:25 PM"; L="0 hours, 15 minutes"}
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("L$r").Value = $row.L
    # Touch the blank STATUS-gap (K) and trailing (M) columns so they
    # materialize as empty cells, matching the existing row layout.
    $ws.Range("K$r").Font.Bold = $false
    $ws.Range("M$r").Font.Bold = $false
}
